$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to match the latest scrape

$ws.Cells.Item(2, 4).Value2 = '70.606.02'
$ws.Cells.Item(2, 5).Value2 = '  +2.05%  '
$ws.Cells.Item(3, 4).Value2 = '3.563.99'
$ws.Cells.Item(3, 5).Value2 = '  +1.27%  '
$ws.Cells.Item(4, 5).Value2 = '  +0.04%  '
$ws.Cells.Item(5, 4).Value2 = '''616.28'
$ws.Cells.Item(5, 5).Value2 = '  +6.44%  '
$ws.Cells.Item(6, 4).Value2 = '''172.23'
$ws.Cells.Item(6, 5).Value2 = '  +0.45%  '
$ws.Cells.Item(7, 4).Value2 = '''0.618'
$ws.Cells.Item(7, 5).Value2 = '  +1.61%  '
$ws.Cells.Item(8, 4).Value2 = '3.559.66'
$ws.Cells.Item(8, 5).Value2 = '  +1.47%  '
$ws.Cells.Item(9, 5).Value2 = '  +0.00%  '
$ws.Cells.Item(10, 4).Value2 = '''0.196'
$ws.Cells.Item(10, 5).Value2 = '  +4.13%  '
$ws.Cells.Item(11, 4).Value2 = '''7.25'
$ws.Cells.Item(11, 5).Value2 = '  +11.05%  '
$ws.Cells.Item(12, 4).Value2 = '''0.586'
$ws.Cells.Item(12, 5).Value2 = '  +0.76%  '
$ws.Cells.Item(13, 4).Value2 = '''46.72'
$ws.Cells.Item(13, 5).Value2 = '  -0.10%  '
$ws.Cells.Item(14, 5).Value2 = '  +1.28%  '
$ws.Cells.Item(15, 4).Value2 = '4.141.61'
$ws.Cells.Item(15, 5).Value2 = '  +1.63%  '
$ws.Cells.Item(16, 4).Value2 = '''8.39'
$ws.Cells.Item(16, 5).Value2 = '  -1.95%  '
$ws.Cells.Item(17, 4).Value2 = '''617.84'
$ws.Cells.Item(17, 5).Value2 = '  -0.71%  '
$ws.Cells.Item(18, 4).Value2 = '3.565.06'
$ws.Cells.Item(18, 5).Value2 = '  +1.65%  '
$ws.Cells.Item(19, 4).Value2 = '70.712.73'
$ws.Cells.Item(19, 5).Value2 = '  +2.32%  '
$ws.Cells.Item(20, 5).Value2 = '  -2.16%  '
$ws.Cells.Item(21, 4).Value2 = '''17.41'
$ws.Cells.Item(21, 5).Value2 = '  -0.31%  '
$ws.Cells.Item(22, 5).Value2 = '  -0.05%  '
$ws.Cells.Item(23, 4).Value2 = '''9.43'
$ws.Cells.Item(23, 5).Value2 = '  -15.40%  '
$ws.Cells.Item(24, 4).Value2 = '''15.77'
$ws.Cells.Item(24, 5).Value2 = '  -1.12%  '
$ws.Cells.Item(25, 4).Value2 = '''96.79'
$ws.Cells.Item(25, 5).Value2 = '  -0.66%  '
$ws.Cells.Item(26, 5).Value2 = '  +1.15%  '
$ws.Cells.Item(27, 5).Value2 = '  +0.02%  '
$ws.Cells.Item(28, 4).Value2 = '''2.61'
$ws.Cells.Item(28, 5).Value2 = '  -0.89%  '
$ws.Cells.Item(29, 4).Value2 = '''33.59'
$ws.Cells.Item(29, 5).Value2 = '  +2.84%  '
$ws.Cells.Item(30, 4).Value2 = '''9.05'
$ws.Cells.Item(30, 5).Value2 = '  -3.05%  '
$ws.Cells.Item(31, 4).Value2 = '''8.51'
$ws.Cells.Item(31, 5).Value2 = '  -0.54%  '
$ws.Cells.Item(32, 4).Value2 = '''3.05'
$ws.Cells.Item(32, 5).Value2 = '  -3.39%  '
$ws.Cells.Item(33, 4).Value2 = '''1.30'
$ws.Cells.Item(33, 5).Value2 = '  -1.05%  '
$ws.Cells.Item(34, 4).Value2 = '''6.95'
$ws.Cells.Item(34, 5).Value2 = '  -0.70%  '
$ws.Cells.Item(35, 4).Value2 = '''574.50'
$ws.Cells.Item(35, 5).Value2 = '  -9.35%  '
$ws.Cells.Item(36, 5).Value2 = '  -1.27%  '
$ws.Cells.Item(37, 4).Value2 = '''3.61'
$ws.Cells.Item(37, 5).Value2 = '  +5.37%  '
$ws.Cells.Item(38, 4).Value2 = '''10.84'
$ws.Cells.Item(38, 5).Value2 = '  +0.89%  '
$ws.Cells.Item(39, 4).Value2 = '''57.49'
$ws.Cells.Item(39, 5).Value2 = '  +1.48%  '
$ws.Cells.Item(40, 4).Value2 = '''0.0471'
$ws.Cells.Item(40, 5).Value2 = '  +5.25%  '
$ws.Cells.Item(41, 5).Value2 = '  +0.30%  '
$ws.Cells.Item(42, 4).Value2 = '''0.142'
$ws.Cells.Item(42, 5).Value2 = '  +4.25%  '
$ws.Cells.Item(43, 4).Value2 = '3.386.80'
$ws.Cells.Item(43, 5).Value2 = '  +0.42%  '
$ws.Cells.Item(44, 4).Value2 = '''0.320'
$ws.Cells.Item(44, 5).Value2 = '  -2.36%  '
$ws.Cells.Item(45, 2).Value2 = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(45, 4).Value2 = '''33.03'
$ws.Cells.Item(45, 5).Value2 = '  +0.31%  '
$ws.Cells.Item(46, 2).Value2 = 'ThetaToken'
$ws.Cells.Item(46, 3).Value2 = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(46, 4).Value2 = '''2.97'
$ws.Cells.Item(46, 5).Value2 = '  +7.08%  '
$ws.Cells.Item(47, 4).Value2 = '0.0₃0701'
$ws.Cells.Item(47, 5).Value2 = '  +1.43%  '
$ws.Cells.Item(48, 4).Value2 = '''2.62'
$ws.Cells.Item(48, 5).Value2 = '  +2.10%  '
$ws.Cells.Item(49, 5).Value2 = '  +0.22%  '
$ws.Cells.Item(50, 4).Value2 = '''133.54'
$ws.Cells.Item(50, 5).Value2 = '  +1.07%  '
$ws.Cells.Item(51, 2).Value2 = 'USDe'
$ws.Cells.Item(51, 3).Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(51, 4).Value2 = '''1.00'
$ws.Cells.Item(51, 5).Value2 = '  -0.01%  '

Write-Output "Applied 96 cell updates"
